$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "Summary": refresh aggregate stats now that trade #100 closed
# and trade #133 was opened.
# -----------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.15    # Current Capital
$summary.Range("B4").Value = 0.94       # Total P&L $
$summary.Range("B5").Value = 0.19       # Total P&L %
$summary.Range("B6").Value = 100        # Total Trades
$summary.Range("B7").Value = 48         # Winning Trades
$summary.Range("B9").Value = 48         # Win Rate %

# -----------------------------------------------------------------
# Sheet "Strategy Status": MarketMaking row (row 5)
# -----------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.15      # Capital
$status.Range("D5").Value = 67          # Trades
$status.Range("E5").Value = 0.83        # P&L $
$status.Range("F5").Value = 1.15        # P&L %
$status.Range("G5").Value = 50.75       # Win Rate %

# -----------------------------------------------------------------
# Sheet "All Trades": close trade #100 (row 101) ...
# -----------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G101").Value = 0.94
$allTrades.Range("H101").Value = "CLOSED"
$allTrades.Range("I101").Value = 4.4444
$allTrades.Range("J101").Value = 0.04
$allTrades.Range("K101").Value = 101.15
$allTrades.Range("L101").Value = "early_exit"
$allTrades.Range("M101").Value = 0.11

# ... and append the newly opened trade #133 as row 134.
$allTrades.Range("A134").Value = 133
$allTrades.Range("B134").NumberFormat = "@"
$allTrades.Range("B134").Value = "2026-02-17"
$allTrades.Range("C134").NumberFormat = "@"
$allTrades.Range("C134").Value = "21:12:38"
$allTrades.Range("D134").Value = "MarketMaking"
$allTrades.Range("E134").Value = "UP"
$allTrades.Range("F134").Value = 0.9
$allTrades.Range("H134").Value = "OPEN"
$allTrades.Range("I134").Value = 0
$allTrades.Range("J134").Value = 0
$allTrades.Range("K134").Value = 101.1096151053151
$allTrades.Range("M134").Value = 0
$allTrades.Range("N134").Value = 0
$allTrades.Range("O134").Value = 0
$allTrades.Range("P134").Value = 0.6
$allTrades.Range("Q134").Value = "Normal spread capture: 19600 bps"

# -----------------------------------------------------------------
# Sheet "MarketMaking": mirror of the same two trades, different
# column layout (L/M/N/O/P/Q order differs from "All Trades").
# -----------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G68").Value = 0.94
$mm.Range("H68").Value = "CLOSED"
$mm.Range("I68").Value = 4.4444
$mm.Range("J68").Value = 0.04
$mm.Range("K68").Value = 101.15
$mm.Range("P68").Value = "early_exit"
$mm.Range("Q68").Value = 0.11

$mm.Range("A101").Value = 133
$mm.Range("B101").NumberFormat = "@"
$mm.Range("B101").Value = "2026-02-17"
$mm.Range("C101").NumberFormat = "@"
$mm.Range("C101").Value = "21:12:38"
$mm.Range("D101").Value = "MarketMaking"
$mm.Range("E101").Value = "UP"
$mm.Range("F101").Value = 0.9
$mm.Range("H101").Value = "OPEN"
$mm.Range("I101").Value = 0
$mm.Range("J101").Value = 0
$mm.Range("K101").Value = 101.1096151053151
$mm.Range("L101").Value = 0
$mm.Range("M101").Value = 0
$mm.Range("N101").Value = 0.6
$mm.Range("O101").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q101").Value = 0
